$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShape = $range.InlineShapes.Item($i)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        [void]$shape.ConvertToInlineShape()
    }
}

# Footer with the first Pearson Edexcel logo (docPr id="3") -> image1.png -> image2.png
$footerA = $sec.Footers.Item(2)
Rename-InlinePicture $footerA.Range "image2.png"

# Footer with the second Pearson Edexcel logo (docPr id="2") -> image1.png -> image2.png
$footerB = $sec.Footers.Item(1)
Rename-InlinePicture $footerB.Range "image2.png"

# Header with the BTec logo (docPr id="1") -> image2.jpg -> image1.jpg
$header = $sec.Headers.Item(2)
Rename-InlinePicture $header.Range "image1.jpg"
